# Update the "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
# timestamps on the zh-cn and de-de report sheets, as produced by a re-run of the
# handback status report generation.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 20:54:08"
$wsZhCn.Range("H2").Value = "2016-03-17 20:54:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 20:54:11"
$wsDeDe.Range("H2").Value = "2016-03-17 20:54:31"
